$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (2..77) down by one row (bottom-up so we
# never overwrite data before it has been read). This makes room for a
# brand new record at row 2 while preserving every other row's values and
# date-column ("s"=2) formatting, which auto-propagates because rows
# 2..77 already share that formatting with their neighbours.
for ($i = 77; $i -ge 2; $i--) {
    $j = $i + 1
    $ws.Range("A$j").Value = $ws.Range("A$i").Value2
    $ws.Range("B$j").Value = $ws.Range("B$i").Value2
    $ws.Range("C$j").Value = $ws.Range("C$i").Value2
    $ws.Range("D$j").Value = $ws.Range("D$i").Value2
    $ws.Range("E$j").Value = $ws.Range("E$i").Value2
    $ws.Range("F$j").Value = $ws.Range("F$i").Value2
    $ws.Range("G$j").Value = $ws.Range("G$i").Value2
    $ws.Range("H$j").Value = $ws.Range("H$i").Value2
    $ws.Range("I$j").Value = $ws.Range("I$i").Value2
    $ws.Range("J$j").Value = $ws.Range("J$i").Value2
    $ws.Range("K$j").Value = $ws.Range("K$i").Value2
    $ws.Range("L$j").Value = $ws.Range("L$i").Value2
    $ws.Range("M$j").Value = $ws.Range("M$i").Value2
    $ws.Range("N$j").Value = $ws.Range("N$i").Value2
    $ws.Range("O$j").Value = $ws.Range("O$i").Value2
    $ws.Range("P$j").Value = $ws.Range("P$i").Value2
    $ws.Range("Q$j").Value = $ws.Range("Q$i").Value2
    $ws.Range("R$j").Value = $ws.Range("R$i").Value2
}

# Row 78 is a brand-new cell for column D (date), so give it the same
# date number format used by every other cell in that column.
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D78").Value = 44489

# Populate the new record at row 2 (weekly price update).
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = 44860
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 300000000
$ws.Range("G2").Value = "Espárragos"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 1000
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 1000
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
